$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add three new header cells (AD1:AF1) with the same formatting as the
# existing header row (bold font, thin border, centered/top aligned).
$headerRng = $ws.Range("AD1:AF1")
$headerRng.Font.Bold = $true
$headerRng.HorizontalAlignment = -4108
$headerRng.VerticalAlignment = -4160
$headerRng.Borders.LineStyle = 1

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Populate the team record (Wins/Losses/Ties) for every data row (2-54).
for ($r = 2; $r -le 54; $r++) {
    $ws.Cells.Item($r, 30).Value = 85
    $ws.Cells.Item($r, 31).Value = 77
    $ws.Cells.Item($r, 32).Value = 0
}
